$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 44449
$ws.Cells.Item(2, 10).Value = 25
$ws.Cells.Item(2, 11).Value = 80000
$ws.Cells.Item(2, 12).Value = 80000
$ws.Cells.Item(2, 13).Value = 80000
$ws.Cells.Item(2, 16).Value = 3200

# Row 3
$ws.Cells.Item(3, 4).Value = 44449
$ws.Cells.Item(3, 9).Value = 'Segunda'
$ws.Cells.Item(3, 10).Value = 20
$ws.Cells.Item(3, 11).Value = 75000
$ws.Cells.Item(3, 12).Value = 75000
$ws.Cells.Item(3, 13).Value = 75000
$ws.Cells.Item(3, 16).Value = 5000

# Row 4
$ws.Cells.Item(4, 4).Value = 44553
$ws.Cells.Item(4, 8).Value = 'Inferno'
$ws.Cells.Item(4, 10).Value = 35
$ws.Cells.Item(4, 11).Value = 45000
$ws.Cells.Item(4, 12).Value = 45000
$ws.Cells.Item(4, 13).Value = 45000
$ws.Cells.Item(4, 16).Value = 1800

# Row 5
$ws.Cells.Item(5, 4).Value = 44425
$ws.Cells.Item(5, 10).Value = 15
$ws.Cells.Item(5, 11).Value = 75000
$ws.Cells.Item(5, 12).Value = 75000
$ws.Cells.Item(5, 13).Value = 75000
$ws.Cells.Item(5, 16).Value = 3000

# Row 6
$ws.Cells.Item(6, 4).Value = 44340
$ws.Cells.Item(6, 8).Value = 'Americana (o)'
$ws.Cells.Item(6, 10).Value = 15

# Row 8
$ws.Cells.Item(8, 4).Value = 44544
$ws.Cells.Item(8, 8).Value = 'Inferno'
$ws.Cells.Item(8, 10).Value = 12
$ws.Cells.Item(8, 11).Value = 35000
$ws.Cells.Item(8, 12).Value = 35000
$ws.Cells.Item(8, 13).Value = 35000
$ws.Cells.Item(8, 16).Value = 1400

# Row 9
$ws.Cells.Item(9, 4).Value = 44460
$ws.Cells.Item(9, 10).Value = 30
$ws.Cells.Item(9, 11).Value = 95000
$ws.Cells.Item(9, 12).Value = 95000
$ws.Cells.Item(9, 13).Value = 95000
$ws.Cells.Item(9, 16).Value = 3800

# Row 10
$ws.Cells.Item(10, 4).Value = 44343
$ws.Cells.Item(10, 9).Value = 'Primera'
$ws.Cells.Item(10, 11).Value = 36000
$ws.Cells.Item(10, 12).Value = 36000
$ws.Cells.Item(10, 13).Value = 36000
$ws.Cells.Item(10, 14).Value = '$/caja 25 kilos'
$ws.Cells.Item(10, 16).Value = 1440
$ws.Cells.Item(10, 17).Value = 25

# Row 13
$ws.Cells.Item(13, 4).Value = 44319
$ws.Cells.Item(13, 8).Value = 'Americana (o)'
$ws.Cells.Item(13, 10).Value = 20
$ws.Cells.Item(13, 11).Value = 30000
$ws.Cells.Item(13, 12).Value = 30000
$ws.Cells.Item(13, 13).Value = 30000
$ws.Cells.Item(13, 16).Value = 1200

# Row 14
$ws.Cells.Item(14, 4).Value = 44193
$ws.Cells.Item(14, 10).Value = 15
$ws.Cells.Item(14, 11).Value = 46000
$ws.Cells.Item(14, 12).Value = 46000
$ws.Cells.Item(14, 13).Value = 46000
$ws.Cells.Item(14, 14).Value = '$/caja 15 kilos'
$ws.Cells.Item(14, 16).Value = 3067
$ws.Cells.Item(14, 17).Value = 15

# Row 15
$ws.Cells.Item(15, 4).Value = 44581
$ws.Cells.Item(15, 9).Value = 'Segunda'
$ws.Cells.Item(15, 10).Value = 30
$ws.Cells.Item(15, 11).Value = 17000
$ws.Cells.Item(15, 12).Value = 17000
$ws.Cells.Item(15, 13).Value = 17000
$ws.Cells.Item(15, 16).Value = 680

# Row 16
$ws.Cells.Item(16, 4).Value = 44474
$ws.Cells.Item(16, 8).Value = 'Americana (o)'
$ws.Cells.Item(16, 10).Value = 18
$ws.Cells.Item(16, 11).Value = 100000
$ws.Cells.Item(16, 12).Value = 100000
$ws.Cells.Item(16, 13).Value = 100000
$ws.Cells.Item(16, 14).Value = '$/caja 25 kilos'
$ws.Cells.Item(16, 16).Value = 4000
$ws.Cells.Item(16, 17).Value = 25

# Row 17
$ws.Cells.Item(17, 4).Value = 44446
$ws.Cells.Item(17, 9).Value = 'Primera'
$ws.Cells.Item(17, 10).Value = 5
$ws.Cells.Item(17, 11).Value = 78000
$ws.Cells.Item(17, 12).Value = 78000
$ws.Cells.Item(17, 13).Value = 78000
$ws.Cells.Item(17, 16).Value = 3120

# Row 18
$ws.Cells.Item(18, 4).Value = 44446
$ws.Cells.Item(18, 8).Value = 'Inferno'
$ws.Cells.Item(18, 10).Value = 4
$ws.Cells.Item(18, 11).Value = 80000
$ws.Cells.Item(18, 12).Value = 80000
$ws.Cells.Item(18, 13).Value = 80000
$ws.Cells.Item(18, 14).Value = '$/caja 15 kilos'
$ws.Cells.Item(18, 16).Value = 5333
$ws.Cells.Item(18, 17).Value = 15
